$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(7).Insert()

$ws.Range("G1").Value = "sum_assured"
$ws.Range("G2:G5").Value = 500000
